# Auto-generated edit script applying the cryptos.xlsx diff.
# Updates the Price (D) and Volume(1h) (E) columns for the
# refreshed coinranking.com snapshot, and swaps the Litecoin/
# ShibaInu (rows 17-18) and Mantle/USDD (rows 50-51) row data
# back to their reported order.
#
# Every touched cell is forced to Text format before the write
# so numeric-looking strings (prices like '93.10' or '1.00',
# percentages, etc.) keep their exact original text formatting
# instead of being auto-coerced to numbers (which would drop
# significant trailing zeros), matching the source workbook
# where every Coin/Link/Price/Volume cell is stored as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.278.36'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.561.18'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.60%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.38'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.95%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -4.35%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.55%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '17.78'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.54%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0783'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.776.53'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -3.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.562.17'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.50%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.56%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.292.66'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.22%  '
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0₃0713'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.13%  '
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '59.28'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.56%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '186.82'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.12'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.27'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.22%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.64%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '140.92'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.97%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.38%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.62%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -4.10%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -6.71%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.41%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.98'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.86%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.088.95'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.34'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.06%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.53%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.15%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.772'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -8.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.797'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +6.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '93.10'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -5.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.11'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.42%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.692.67'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₆0112'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.23%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '52.48'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.85%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.07%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.29%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.404'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.84%  '
